# "Generate Report for Archive"
# The localization status report is regenerated: the "Ready for handoff"
# status is updated to "In Translation" everywhere it appears (Overview
# summary sheet + the per-locale detail sheets), and the now-narrower
# "Status"/locale-status columns are resized to fit the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text wherever it appears ---

# Overview sheet: per-locale status columns (E = zh-cn, F = de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Per-locale detail sheets: Status column (C)
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Resize the affected columns now that the text is shorter ---
# (stored column width in the sheet XML is ColumnWidth + 5/6 "characters";
# this expresses the new narrower width the report generator produced)
$newColumnWidth = 12.576851254417766

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth

$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
